$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017198638231408
$ws.Range("D2").Value = 1.019259137241189
$ws.Range("E2").Value = 1.018627844743263
$ws.Range("F2").Value = 1.015556486117394
$ws.Range("I2").Value = 1.027201412240472
$ws.Range("J2").Value = 1.022414116398068
$ws.Range("K2").Value = 1.022102915582371
$ws.Range("L2").Value = 1.021473496824078
$ws.Range("M2").Value = 1.018411290498982
$ws.Range("N2").Value = 1.011565174197355
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018087425805349
$ws.Range("D3").Value = 1.020017845863896
$ws.Range("E3").Value = 1.019379810814981
$ws.Range("F3").Value = 1.01708642831064
$ws.Range("I3").Value = 1.027276365072362
$ws.Range("J3").Value = 1.022939005637651
$ws.Range("K3").Value = 1.022667740814182
$ws.Range("L3").Value = 1.022031462339788
$ws.Range("M3").Value = 1.019744413467536
$ws.Range("N3").Value = 1.011740392504969
$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.018662635385054
$ws.Range("D4").Value = 1.020509183122047
$ws.Range("E4").Value = 1.019866869715271
$ws.Range("F4").Value = 1.01807644942963
$ws.Range("I4").Value = 1.027323310069579
$ws.Range("J4").Value = 1.023278145264374
$ws.Range("K4").Value = 1.023032964237131
$ws.Range("L4").Value = 1.022392329977132
$ws.Range("M4").Value = 1.020606602318464
$ws.Range("N4").Value = 1.011853555426933
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018904478093647
$ws.Range("D5").Value = 1.020715836732619
$ws.Range("E5").Value = 1.020071744883341
$ws.Range("F5").Value = 1.018492670002252
$ws.Range("I5").Value = 1.02734267321877
$ws.Range("J5").Value = 1.023420599768984
$ws.Range("K5").Value = 1.023186442311653
$ws.Range("L5").Value = 1.022543996688174
$ws.Range("M5").Value = 1.020968967244689
$ws.Range("N5").Value = 1.01190107745994
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018945085978464
$ws.Range("D6").Value = 1.020750540355014
$ws.Range("E6").Value = 1.0201061510791
$ws.Range("F6").Value = 1.018562556322634
$ws.Range("I6").Value = 1.027345902520503
$ws.Range("J6").Value = 1.023444511471256
$ws.Range("K6").Value = 1.023212208314003
$ws.Range("L6").Value = 1.022569459715062
$ws.Range("M6").Value = 1.02102980420192
$ws.Range("N6").Value = 1.011909053583837
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018665866803051
$ws.Range("D7").Value = 1.02051194406379
$ws.Range("E7").Value = 1.019869606814159
$ws.Range("F7").Value = 1.018082010921866
$ws.Range("I7").Value = 1.027323570265181
$ws.Range("J7").Value = 1.023280049219634
$ws.Range("K7").Value = 1.023035015262353
$ws.Range("L7").Value = 1.022394356721042
$ws.Range("M7").Value = 1.020611444641575
$ws.Range("N7").Value = 1.011854190621878
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017498986235885
$ws.Range("D8").Value = 1.019515462633087
$ws.Range("E8").Value = 1.018881873923683
$ws.Range("F8").Value = 1.016073530463887
$ws.Range("I8").Value = 1.027227064543197
$ws.Range("J8").Value = 1.022591608229865
$ws.Range("K8").Value = 1.022293853677009
$ws.Range("L8").Value = 1.021662099454596
$ws.Range("M8").Value = 1.018861917195906
$ws.Range("N8").Value = 1.011624434463434
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.015443613376349
$ws.Range("D9").Value = 1.017762650095319
$ws.Range("E9").Value = 1.01714512727862
$ws.Range("F9").Value = 1.012534477933386
$ws.Range("I9").Value = 1.027045120060954
$ws.Range("J9").Value = 1.021374683196139
$ws.Range("K9").Value = 1.020985891190543
$ws.Range("L9").Value = 1.020370464288221
$ws.Range("M9").Value = 1.015775552657288
$ws.Range("N9").Value = 1.011217935524287
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.014073938767411
$ws.Range("D10").Value = 1.016596249396817
$ws.Range("E10").Value = 1.015989877589425
$ws.Range("F10").Value = 1.010174917073378
$ws.Range("I10").Value = 1.026915851886397
$ws.Range("J10").Value = 1.020560862302036
$ws.Range("K10").Value = 1.020112638614536
$ws.Range("L10").Value = 1.019508527237093
$ws.Range("M10").Value = 1.013715396330912
$ws.Range("N10").Value = 1.010945843915559
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013480994664971
$ws.Range("D11").Value = 1.016091702103022
$ws.Range("E11").Value = 1.01549026404492
$ws.Range("F11").Value = 1.009153090336653
$ws.Range("I11").Value = 1.026857990890147
$ws.Range("J11").Value = 1.020207870652282
$ws.Range("K11").Value = 1.019734212828705
$ws.Range("L11").Value = 1.019135103908121
$ws.Range("M11").Value = 1.012822659103832
$ws.Range("N11").Value = 1.010827768040863
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.013260769189277
$ws.Range("D12").Value = 1.015904368370464
$ws.Range("E12").Value = 1.015304778813519
$ws.Range("F12").Value = 1.008773514483149
$ws.Range("I12").Value = 1.026836215547974
$ws.Range("J12").Value = 1.020076663641128
$ws.Range("K12").Value = 1.019593603680047
$ws.Range("L12").Value = 1.018996368383926
$ws.Range("M12").Value = 1.012490950395568
$ws.Range("N12").Value = 1.010783870708023
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.013308007376244
$ws.Range("D13").Value = 1.015944548569586
$ws.Range("E13").Value = 1.015344561782846
$ws.Range("F13").Value = 1.008854935965481
$ws.Range("I13").Value = 1.02684089924283
$ws.Range("J13").Value = 1.020104812068426
$ws.Range("K13").Value = 1.019623766848696
$ws.Range("L13").Value = 1.01902612895203
$ws.Range("M13").Value = 1.012562107929457
$ws.Range("N13").Value = 1.010793288583959
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01346279034172
$ws.Range("D14").Value = 1.016076215451329
$ws.Range("E14").Value = 1.015474929866423
$ws.Range("F14").Value = 1.009121714998813
$ws.Range("I14").Value = 1.026856196708801
$ws.Range("J14").Value = 1.020197026878501
$ws.Range("K14").Value = 1.019722590947783
$ws.Range("L14").Value = 1.019123636585311
$ws.Range("M14").Value = 1.012795242155926
$ws.Range("N14").Value = 1.010824140262285
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013558159940474
$ws.Range("D15").Value = 1.016157350065615
$ws.Range("E15").Value = 1.015555266355525
$ws.Range("F15").Value = 1.009286083027923
$ws.Range("I15").Value = 1.026865584466091
$ws.Range("J15").Value = 1.020253831521261
$ws.Range("K15").Value = 1.0197834737895
$ws.Range("L15").Value = 1.01918371036299
$ws.Range("M15").Value = 1.012938869610557
$ws.Range("N15").Value = 1.010843143875774
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.014113293358209
$ws.Range("D16").Value = 1.016629745383291
$ws.Range("E16").Value = 1.016023048363371
$ws.Range("F16").Value = 1.010242729308241
$ws.Range("I16").Value = 1.026919652199422
$ws.Range("J16").Value = 1.020584276534528
$ws.Range("K16").Value = 1.020137747162761
$ws.Range("L16").Value = 1.019533305958094
$ws.Range("M16").Value = 1.013774629709844
$ws.Range("N16").Value = 1.010953674789004
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014461549932891
$ws.Range("D17").Value = 1.01692620402631
$ws.Range("E17").Value = 1.016316641436165
$ws.Range("F17").Value = 1.010842771764093
$ws.Range("I17").Value = 1.026953062471262
$ws.Range("J17").Value = 1.020791394992847
$ws.Range("K17").Value = 1.020359892982586
$ws.Range("L17").Value = 1.019752545078804
$ws.Range("M17").Value = 1.014298695588502
$ws.Range("N17").Value = 1.011022938817354
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014664695018996
$ws.Range("D18").Value = 1.017099172682769
$ws.Range("E18").Value = 1.016487948808224
$ws.Range("F18").Value = 1.011192755053749
$ws.Range("I18").Value = 1.026972368026338
$ws.Range("J18").Value = 1.020912145602683
$ws.Range("K18").Value = 1.02048943775157
$ws.Range("L18").Value = 1.01988040433081
$ws.Range("M18").Value = 1.014604309453199
$ws.Range("N18").Value = 1.011063314411713
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.0147339644278
$ws.Range("D19").Value = 1.017158158873509
$ws.Range("E19").Value = 1.016546370264165
$ws.Range("F19").Value = 1.011312088563091
$ws.Range("I19").Value = 1.026978919826362
$ws.Range("J19").Value = 1.020953308581456
$ws.Range("K19").Value = 1.020533604218233
$ws.Range("L19").Value = 1.019923997759349
$ws.Range("M19").Value = 1.014708505021613
$ws.Range("N19").Value = 1.01107707721155
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014424183960361
$ws.Range("D20").Value = 1.016894391710452
$ws.Range("E20").Value = 1.016285135518084
$ws.Range("F20").Value = 1.010778394098917
$ws.Range("I20").Value = 1.026949496697554
$ws.Range("J20").Value = 1.020769179148488
$ws.Range("K20").Value = 1.020336061848794
$ws.Range("L20").Value = 1.019729024780772
$ws.Range("M20").Value = 1.014242475026038
$ws.Range("N20").Value = 1.011015510021485
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013417210052284
$ws.Range("D21").Value = 1.016037440688802
$ws.Range("E21").Value = 1.015436537128475
$ws.Range("F21").Value = 1.009043155928453
$ws.Range("I21").Value = 1.026851699801339
$ws.Range("J21").Value = 1.020169874403401
$ws.Range("K21").Value = 1.019693490952062
$ws.Range("L21").Value = 1.019094923827971
$ws.Range("M21").Value = 1.012726592911449
$ws.Range("N21").Value = 1.010815056279149
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.012784203523221
$ws.Range("D22").Value = 1.015499090859211
$ws.Range("E22").Value = 1.014903530881369
$ws.Range("F22").Value = 1.0079520016751
$ws.Range("I22").Value = 1.02678857270313
$ws.Range("J22").Value = 1.019792545986983
$ws.Range("K22").Value = 1.019289221160267
$ws.Range("L22").Value = 1.018696069067279
$ws.Range("M22").Value = 1.011772879945684
$ws.Range("N22").Value = 1.010688799235961
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01311976100789
$ws.Range("D23").Value = 1.015784437474132
$ws.Range("E23").Value = 1.015186036023295
$ws.Range("F23").Value = 1.008530458471884
$ws.Range("I23").Value = 1.02682219276131
$ws.Range("J23").Value = 1.019992624301479
$ws.Range("K23").Value = 1.01950355674962
$ws.Range("L23").Value = 1.018907525532283
$ws.Range("M23").Value = 1.012278521393974
$ws.Range("N23").Value = 1.010755751656801
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014441067994937
$ws.Range("D24").Value = 1.016908766174174
$ws.Range("E24").Value = 1.016299371502481
$ws.Range("F24").Value = 1.010807483628345
$ws.Range("I24").Value = 1.026951108479844
$ws.Range("J24").Value = 1.020779217711683
$ws.Range("K24").Value = 1.020346830202317
$ws.Range("L24").Value = 1.019739652650383
$ws.Range("M24").Value = 1.014267878877581
$ws.Range("N24").Value = 1.01101886685162
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.015974877015301
$ws.Range("D25").Value = 1.018215420171924
$ws.Range("E25").Value = 1.017593666566838
$ws.Range("F25").Value = 1.013449420386482
$ws.Range("I25").Value = 1.02709356374465
$ws.Range("J25").Value = 1.021689736196166
$ws.Range("K25").Value = 1.021324258019106
$ws.Range("L25").Value = 1.020704534393467
$ws.Range("M25").Value = 1.016573888643944
$ws.Range("N25").Value = 1.011323218482173
